# Update "想去人数" (interested-count) values on three sheets to reflect
# regenerated output data (gh-pages rebuild at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions) sheet
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 158
$wsExhibit.Range("F9").Value = 273

# 演出 (Performances) sheet
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 84
$wsShow.Range("F3").Value = 35

# 全部类型 (All types, combined) sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 158
$wsAll.Range("F3").Value = 84
$wsAll.Range("F10").Value = 273
$wsAll.Range("F11").Value = 35
